# Weekly driver report update for 2025-04-20
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# "Bad Drivers" table (rows 3-5)
# ---------------------------------------------------------------------------

# Row 3: Intel(R) Wi-Fi 6E AX211 160MHz - 22.230.0.8
$ws.Range("D3").Value = 94.90000000000001

# Row 4: Intel(R) Wi-Fi 6 AX201 160MHz - 23.30.0.6
$ws.Range("B4").Value = 59
$ws.Range("C4").Value = 1684
$ws.Range("D4").Value = 98.3

# Row 5: Totals
$ws.Range("B5").Value = 62
$ws.Range("C5").Value = 1909

# ---------------------------------------------------------------------------
# "Good Drivers (Roaming > 99.8%)" table (rows 13-21)
# A newer driver (23.100.0.4) moved to the top of the vintage-sorted list and
# every other row shifted down one slot; the oldest previous row (21.70.0.6 /
# 113652 / 2019-12-14) dropped out of the visible window.
# ---------------------------------------------------------------------------

# Row 13: Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4
$ws.Range("A13").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4"
$ws.Range("B13").Value = 445055
$ws.Range("D13").Value = 99.90000000000001
$ws.Range("E13").Value = "'2024-11-10"

# Row 14: Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.3.1
$ws.Range("A14").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.3.1"
$ws.Range("B14").Value = 10661
$ws.Range("D14").Value = 100
$ws.Range("E14").Value = "'2022-08-29"

# Row 15: Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.0.3
$ws.Range("A15").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.0.3"
$ws.Range("B15").Value = 14239
$ws.Range("D15").Value = 100
$ws.Range("E15").Value = "'2022-05-23"

# Row 16: Intel(R) Wi-Fi 6E AX211 160MHz - 22.100.1.1
$ws.Range("A16").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.100.1.1"
$ws.Range("B16").Value = 265400
$ws.Range("D16").Value = 99.90000000000001
$ws.Range("E16").Value = "'2022-05-01"

# Row 17: Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9
$ws.Range("A17").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9"
$ws.Range("B17").Value = 77849
$ws.Range("D17").Value = 99.90000000000001
$ws.Range("E17").Value = "'2021-08-18"

# Row 18: Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1
$ws.Range("A18").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1"
$ws.Range("B18").Value = 34244
$ws.Range("D18").Value = 100
$ws.Range("E18").Value = "'2021-04-27"

# Row 19: Intel(R) Wi-Fi 6 AX201 160MHz - 21.110.3.2
$ws.Range("A19").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.110.3.2"
$ws.Range("B19").Value = 59673
$ws.Range("D19").Value = 100
$ws.Range("E19").Value = "'2020-08-05"

# Row 20: Intel(R) Wi-Fi 6 AX201 160MHz - 21.70.0.6
$ws.Range("A20").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.70.0.6"
$ws.Range("B20").Value = 113652
$ws.Range("D20").Value = 100
$ws.Range("E20").Value = "'2020-01-06"

# Row 21: Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1
$ws.Range("A21").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1"
$ws.Range("B21").Value = 56018
$ws.Range("D21").Value = 100
$ws.Range("E21").Value = "'2019-12-14"
